# Generate Report for Handback
#
# The prior handback run included two source files:
#   153d35a8-dc9a-4c9c-b693-8e20c14b471b.md
#   f4b530c9-fe30-4fb0-b7a2-0e6ff5b6abd3.md
# This run only has the first file to report on, and it was handed off/back
# at a later time than the previous run. So for every worksheet we drop the
# row describing the second (f4b530c9...) file, and refresh the
# handoff/handback timestamps recorded for the remaining file.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: drop the row for the f4b530c9... file (row 3) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Rows.Item(3).Delete()

# --- zh-cn sheet: drop its row 3, refresh the handoff/handback datetimes ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Rows.Item(3).Delete()
$wsZhCn.Range("E2").Value = "2016-03-17 16:38:26"
$wsZhCn.Range("H2").Value = "2016-03-17 16:38:44"

# --- de-de sheet: drop its row 3, refresh the handoff/handback datetimes ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Rows.Item(3).Delete()
$wsDeDe.Range("E2").Value = "2016-03-17 16:38:29"
$wsDeDe.Range("H2").Value = "2016-03-17 16:38:49"
